# "CERRADURA PARA MUEBLES.xlsx" price-list refresh:
#   - bump the price-list date in A1 by one day
#   - update the three unit prices in column D (rows 31-33) to the
#     corrected values pulled from the source spreadsheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 45311

$ws.Range("D31").Value = 2022.982
$ws.Range("D32").Value = 835.851
$ws.Range("D33").Value = 1715.33
